$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.733.65"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.696.81"
$ws.Range("E3").Value = "  +2.62%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.04"
$ws.Range("E5").Value = "  +1.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.57"
$ws.Range("E6").Value = "  +1.44%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("E9").Value = "  +4.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.99"
$ws.Range("E10").Value = "  +3.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.36"
$ws.Range("E13").Value = "  +4.36%  "

$ws.Range("E14").Value = "  +8.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.184.92"
$ws.Range("E15").Value = "  +2.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.623.25"
$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.687.64"
$ws.Range("E17").Value = "  +1.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.68"
$ws.Range("E18").Value = "  +1.30%  "

$ws.Range("E19").Value = "  +2.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "358.65"
$ws.Range("E20").Value = "  +1.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.63"
$ws.Range("E21").Value = "  +3.82%  "

$ws.Range("E23").Value = "  +4.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.87"
$ws.Range("E24").Value = "  +3.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000106"
$ws.Range("E25").Value = "  +11.50%  "

$ws.Range("E26").Value = "  +0.68%  "

$ws.Range("E27").Value = "  +2.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.53"
$ws.Range("E28").Value = "  +5.88%  "

$ws.Range("E29").Value = "  +4.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "547.15"
$ws.Range("E30").Value = "  +3.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.20"
$ws.Range("E31").Value = "  +4.90%  "

$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.81"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("E34").Value = "  +6.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("E35").Value = "  -1.82%  "

$ws.Range("E36").Value = "  +1.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.85"
$ws.Range("E37").Value = "  +2.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.93"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "171.55"
$ws.Range("E41").Value = "  +3.92%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.59"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("E44").Value = "  +2.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0617"
$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.60"
$ws.Range("E46").Value = "  +2.47%  "

$ws.Range("E47").Value = "  +2.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0266"
$ws.Range("E48").Value = "  +4.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.655"
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.04"
$ws.Range("E50").Value = "  +8.54%  "

$ws.Range("E51").Value = "  +1.48%  "
